# Add 2022-Q4 data:
#  1. Insert a new worksheet "2022-Q4" right after "总计" (before "2022-Q2").
#  2. Fill it with the quarterly fund-holding table, matching the look of the
#     other quarterly sheets (bold/centered/bordered header row + index col).
#  3. Prepend the new quarter's summary row into "总计" (shifting the rest
#     down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new sheet in the correct tab position.
# ---------------------------------------------------------------------------
$insertBefore = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($insertBefore)
$newSheet.Name = "2022-Q4"

# Re-fetch everything we need *after* the sheet collection changed shape -
# stale references silently no-op on PasteSpecial.
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item("2022-Q4")
$styleSource = $wb.Worksheets.Item("2022-Q2")

# Match the sheet-level look (outline direction, page margins) used by
# every other quarterly sheet.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q4" with the fund table.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0, "013184", "广发恒阳一年持有期混合A", "3.74", "28.15", "0.79", "0.0295", 5),
    @(1, "013185", "广发恒阳一年持有期混合C", "3.21", "28.15", "0.79", "0.0254", 5),
    @(2, "002135", "广发鑫源灵活配置混合A", "0.59", "27.74", "1.63", "0.0096", 2),
    @(3, "090011", "大成核心双动力混合",    "0.26", "93.43", "2.37", "0.0062", 9),
    @(4, "002136", "广发鑫源灵活配置混合C", "0.19", "27.74", "1.63", "0.0031", 2)
)

$rowIdx = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($rowIdx, 1).Value = $row[0]

    # B:G are text columns on every quarterly sheet (fund code keeps leading
    # zeros, the numeric-looking figures stay text) - force text typing so
    # Excel doesn't auto-coerce them to numbers.
    $textRange = $newSheet.Range($newSheet.Cells.Item($rowIdx, 2), $newSheet.Cells.Item($rowIdx, 7))
    $textRange.NumberFormat = "@"
    $newSheet.Cells.Item($rowIdx, 2).Value = $row[1]
    $newSheet.Cells.Item($rowIdx, 3).Value = $row[2]
    $newSheet.Cells.Item($rowIdx, 4).Value = $row[3]
    $newSheet.Cells.Item($rowIdx, 5).Value = $row[4]
    $newSheet.Cells.Item($rowIdx, 6).Value = $row[5]
    $newSheet.Cells.Item($rowIdx, 7).Value = $row[6]
    $textRange.Style = "Normal"

    $newSheet.Cells.Item($rowIdx, 8).Value = $row[7]

    $rowIdx++
}

# Match the bold/centered/bordered look used for the header row and the
# index column on every other quarterly sheet.
$styleSource.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$styleSource.Range("A2:A6").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Update "总计": prepend the 2022-Q4 summary row, pushing the rest down.
# ---------------------------------------------------------------------------
$summary = @(
    @("2022-Q4", 5, 0.07000000000000001),
    @("2022-Q2", 12, 1.31),
    @("2022-Q1", 7, 0.25),
    @("2021-Q2", 1, 0),
    @("2021-Q1", 6, 0.46),
    @("2020-Q4", 8, 0.59)
)

for ($i = 0; $i -lt $summary.Length; $i++) {
    $r = $i + 2
    $entry = $summary[$i]
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $entry[0]
    $totalSheet.Cells.Item($r, 3).Value = $entry[1]
    $totalSheet.Cells.Item($r, 4).Value = $entry[2]
}

# Row 7 is brand new - give its index cell (A7) the same style as the other
# index cells in column A (bold, centered, bordered).
$totalSheet.Cells.Item(2, 1).Copy()
$totalSheet.Cells.Item(7, 1).PasteSpecial(-4122)

# Keep "总计" the active/selected sheet, same as before the edit.
$totalSheet.Activate()
[void]$totalSheet.Range("A1").Select()
